# Staging.OrganizationPersonRole.xlsx template update.
#
# The "SourceKey" columns are being renamed to "BusinessKey" columns, and a
# new "OrganizationBusinessKey" header is introduced ahead of the existing
# "OrganizationPersonRole_ID" header (the latter shifts one column to the
# right). Net effect on row 2 (the example/header row), by cell position:
#   A2: OrganizationPersonRole_ID -> OrganizationBusinessKey
#   B2: OrganizationSourceKey     -> OrganizationPersonRole_ID
#   C2: PersonSourceKey           -> PersonBusinessKey
#   D2: RoleSourceKey             -> RoleBusinessKey
# Row 1 (the "For internal use only..." banner) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "OrganizationBusinessKey"
$ws.Range("B2").Value = "OrganizationPersonRole_ID"
$ws.Range("C2").Value = "PersonBusinessKey"
$ws.Range("D2").Value = "RoleBusinessKey"
